# Q1 sheet: update the goal-seek target and the resulting force-of-interest
# rate, then rework the D (per-period growth factor) and E (accumulation
# factor) formulas.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Q1")

# --- Inputs -----------------------------------------------------------
# L1 = "GOALSEEKED AMOUNT" target, raised from 80,000 to 100,000.
$ws1.Range("L1").Value2 = 100000
# G1 = force-of-interest rate solved by Goal Seek so that G6 (the total
# accumulated value) equals L1.
$ws1.Range("G1").Value2 = 0.069389487671259062

# --- D column: per-period growth factor --------------------------------
# Breakpoints shift one year earlier, and the very last period (B = $I$2)
# now returns 1 instead of falling through to EXP($I$1).
$ws1.Range("D6").Formula = '=IF(B6<=($G$2-1),(1-$G$1/$G$4)^(-$G$4),IF(B6<=($H$2-1),(1+$H$1/$H$4)^($H$4),IF(B6=$I$2,1,EXP($I$1))))'
for ($r = 7; $r -le 36; $r++) {
    $f = '=IF(B' + $r + '<=($G$2-1),(1-$G$1/$G$4)^(-$G$4),IF(B' + $r + '<=($H$2-1),(1+$H$1/$H$4)^($H$4),IF(B' + $r + '=$I$2,1,EXP($I$1))))'
    $ws1.Range("D$r").Formula = $f
}

# --- E column: accumulation factor, rewritten as a running PRODUCT -----
for ($r = 6; $r -le 36; $r++) {
    $f = '=PRODUCT(D' + $r + ':$D$36)'
    $ws1.Range("E$r").Formula = $f
}

# --- C21: re-entered explicitly (same logic, parenthesised reference) --
$ws1.Range("C21").Formula = '=IF(B21<=($C$2),$C$1,$D$1)'

# --- Selection / active sheet bookkeeping ------------------------------
# The user ends the session on Q1 (previously Q2 was the active tab),
# with the cursor on L2 (the goal-seeked rate's label cell next door).
$ws1.Activate()
$ws1.Range("L2").Select()
